# -----------------------------------------------------------------------
# LoginData.xlsx edit: add a new "Sheet2" with customer-registration test
# data, re-point the Sheet1 selection, and make Sheet2 the active sheet.
# -----------------------------------------------------------------------

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet1: change the current selection (no longer the active tab) ---
$ws1.Range("A2:B2").Select()

# --- Insert a brand-new worksheet right after Sheet1 -------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

# --- Header row (row 1) --------------------------------------------------
$ws2.Range("A1").Value = "CustomerLoginName"
$ws2.Range("B1").Value = "CustomerLoginPassword"
$ws2.Range("C1").Value = "CustomerName"
$ws2.Range("D1").Value = "CustomerGender"
$ws2.Range("E1").Value = "CustomerDateOfBirth"
$ws2.Range("F1").Value = "CustomerAddress"
$ws2.Range("G1").Value = "CustomerCity"
$ws2.Range("H1").Value = "CustomerState"
$ws2.Range("I1").Value = "CustomerPin"
$ws2.Range("J1").Value = "CustomerTelephone"
$ws2.Range("K1").Value = "CustomerEmail"
$ws2.Range("L1").Value = "CustomerPassword"

# --- Row 2: put the genuinely-numeric values in first (while the cells --
#     still carry the default General format) so they persist as real
#     numbers once the row is switched to Text format below. -------------
$ws2.Range("E2").Value = 10152022
$ws2.Range("I2").Value = 785680
$ws2.Range("J2").Value = 6541215
$ws2.Range("L2").Value = 1254

$ws2.Range("A2").Value = "mngr435533"
$ws2.Range("B2").Value = "EgUhAte"
$ws2.Range("C2").Value = "Amir Ansari"
$ws2.Range("D2").Value = "male"
$ws2.Range("F2").Value = "INDIA"
$ws2.Range("G2").Value = "Amguri"
$ws2.Range("H2").Value = "Assam"
$ws2.Range("K2").Value = "a23432www4@exy.com"

# --- Row 3: every value here is stored as text (even the numeric-looking
#     ones), so switch the row to Text format before typing anything. ----
$ws2.Range("A3:L3").NumberFormat = "@"
$ws2.Range("A3").Value = "mngr435533"
$ws2.Range("B3").Value = "EgUhAte"
$ws2.Range("C3").Value = "ABCD EFGH"
$ws2.Range("D3").Value = "female"
$ws2.Range("E3").Value = "03021965"
$ws2.Range("F3").Value = "INDIA"
$ws2.Range("G3").Value = "Sivasagar"
$ws2.Range("H3").Value = "Assam"
$ws2.Range("I3").Value = "785681"
$ws2.Range("J3").Value = "65410000"
$ws2.Range("K3").Value = "aasdfdsf@exy.com"
$ws2.Range("L3").Value = "2000"

# --- Now switch row 2 to Text format too (values are already entered, so
#     the real numbers stay numbers; only the display format changes). --
$ws2.Range("A2:L2").NumberFormat = "@"

# --- Pre-format the remaining (still-empty) data rows as Text ----------
$ws2.Range("C4:L13").NumberFormat = "@"

# --- Hyperlinks on the two e-mail cells ---------------------------------
$ws2.Hyperlinks.Add($ws2.Range("K2"), "mailto:a23432www4@exy.com")
$ws2.Hyperlinks.Add($ws2.Range("K3"), "mailto:aasdfdsf@exy.com")

# --- Header formatting: bold + yellow fill ------------------------------
$headerRange = $ws2.Range("A1:L1")
$headerRange.Font.Bold = $true
$headerRange.Interior.Color = 65535

# --- Column widths (approximate best-fit widths from the source file) --
$ws2.Columns.Item(1).ColumnWidth = 17.39
$ws2.Columns.Item(3).ColumnWidth = 12.94
$ws2.Columns.Item(4).ColumnWidth = 13.94
$ws2.Columns.Item(5).ColumnWidth = 17.94
$ws2.Columns.Item(6).ColumnWidth = 14.39
$ws2.Columns.Item(7).ColumnWidth = 11.17
$ws2.Columns.Item(8).ColumnWidth = 12.28
$ws2.Columns.Item(10).ColumnWidth = 16.72
$ws2.Columns.Item(11).ColumnWidth = 24.28
$ws2.Columns.Item(12).ColumnWidth = 15.94

# --- Page setup (so a <pageSetup> element is emitted, like Sheet1) -----
$ws2.PageSetup.Orientation = 1

# --- Selection / active sheet -------------------------------------------
$ws2.Range("K3").Select()
$ws2.Activate()

Write-Host "Sheet2 added with customer registration test data."
